$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.676342725753784
$ws.Range("B1").Value = 4.045720100402832
$ws.Range("C1").Value = 3.305055856704712
$ws.Range("D1").Value = 1.604742527008057
$ws.Range("E1").Value = 0.7457998394966125
